$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fases")
Write-Host $ws.Name
